$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-10: new TPM-derived values after adding "Resolving-Mac" cluster
# Row 2
$ws.Cells.Item(2, 7).Value = 3.295960333333333
$ws.Cells.Item(2, 8).Value = 9.887881
$ws.Cells.Item(2, 9).Value = 0.02344673072001071
$ws.Cells.Item(2, 10).Value = 0.02344673072001072
$ws.Cells.Item(2, 13).Value = 22.618885
$ws.Cells.Item(2, 14).Value = 67.856655
$ws.Cells.Item(2, 15).Value = 0.9027998993061069
$ws.Cells.Item(2, 16).Value = 0.902799899306107
$ws.Cells.Item(2, 17).Value = 74.55094774422834
$ws.Cells.Item(2, 18).Value = 670.958529698055
$ws.Cells.Item(2, 19).Value = 0.02116770613308308
$ws.Cells.Item(2, 20).Value = 0.02116770613308308

# Row 3
$ws.Cells.Item(3, 7).Value = 3.295960333333333
$ws.Cells.Item(3, 8).Value = 9.887881
$ws.Cells.Item(3, 9).Value = 0.02344673072001071
$ws.Cells.Item(3, 10).Value = 0.02344673072001072
$ws.Cells.Item(3, 13).Value = 2.154739666666666
$ws.Cells.Item(3, 14).Value = 6.464219
$ws.Cells.Item(3, 15).Value = 0.08600330007856447
$ws.Cells.Item(3, 16).Value = 0.08600330007856449
$ws.Cells.Item(3, 17).Value = 7.101936469993221
$ws.Cells.Item(3, 18).Value = 63.917428229939
$ws.Cells.Item(3, 19).Value = 0.002016496217974377
$ws.Cells.Item(3, 20).Value = 0.002016496217974378

# Row 4
$ws.Cells.Item(4, 7).Value = 3.295960333333333
$ws.Cells.Item(4, 8).Value = 9.887881
$ws.Cells.Item(4, 9).Value = 0.02344673072001071
$ws.Cells.Item(4, 10).Value = 0.02344673072001072
$ws.Cells.Item(4, 13).Value = 0.2805263333333333
$ws.Cells.Item(4, 14).Value = 0.841579
$ws.Cells.Item(4, 15).Value = 0.0111968006153285
$ws.Cells.Item(4, 16).Value = 0.01119680061532851
$ws.Cells.Item(4, 17).Value = 0.9246036671221111
$ws.Cells.Item(4, 18).Value = 8.321433004099
$ws.Cells.Item(4, 19).Value = 0.0002625283689532577
$ws.Cells.Item(4, 20).Value = 0.0002625283689532578

# Row 5
$ws.Cells.Item(5, 7).Value = 136.9593173333334
$ws.Cells.Item(5, 8).Value = 410.8779520000001
$ws.Cells.Item(5, 9).Value = 0.9742982039663998
$ws.Cells.Item(5, 10).Value = 0.9742982039663999
$ws.Cells.Item(5, 13).Value = 22.618885
$ws.Cells.Item(5, 14).Value = 67.856655
$ws.Cells.Item(5, 15).Value = 0.9027998993061069
$ws.Cells.Item(5, 16).Value = 0.902799899306107
$ws.Cells.Item(5, 17).Value = 3097.867048441174
$ws.Cells.Item(5, 18).Value = 27880.80343597056
$ws.Cells.Item(5, 19).Value = 0.8795963204349866
$ws.Cells.Item(5, 20).Value = 0.8795963204349867

# Row 6
$ws.Cells.Item(6, 7).Value = 136.9593173333334
$ws.Cells.Item(6, 8).Value = 410.8779520000001
$ws.Cells.Item(6, 9).Value = 0.9742982039663998
$ws.Cells.Item(6, 10).Value = 0.9742982039663999
$ws.Cells.Item(6, 13).Value = 2.154739666666666
$ws.Cells.Item(6, 14).Value = 6.464219
$ws.Cells.Item(6, 15).Value = 0.08600330007856447
$ws.Cells.Item(6, 16).Value = 0.08600330007856449
$ws.Cells.Item(6, 17).Value = 295.1116737777209
$ws.Cells.Item(6, 18).Value = 2656.005063999488
$ws.Cells.Item(6, 19).Value = 0.08379286080172868
$ws.Cells.Item(6, 20).Value = 0.08379286080172872

# Row 7
$ws.Cells.Item(7, 7).Value = 136.9593173333334
$ws.Cells.Item(7, 8).Value = 410.8779520000001
$ws.Cells.Item(7, 9).Value = 0.9742982039663998
$ws.Cells.Item(7, 10).Value = 0.9742982039663999
$ws.Cells.Item(7, 13).Value = 0.2805263333333333
$ws.Cells.Item(7, 14).Value = 0.841579
$ws.Cells.Item(7, 15).Value = 0.0111968006153285
$ws.Cells.Item(7, 16).Value = 0.01119680061532851
$ws.Cells.Item(7, 17).Value = 38.42069510735645
$ws.Cells.Item(7, 18).Value = 345.786255966208
$ws.Cells.Item(7, 19).Value = 0.01090902272968444
$ws.Cells.Item(7, 20).Value = 0.01090902272968445

# Row 8
$ws.Cells.Item(8, 7).Value = 0.1784693333333333
$ws.Cells.Item(8, 8).Value = 0.535408
$ws.Cells.Item(8, 9).Value = 0.001269591250272884
$ws.Cells.Item(8, 10).Value = 0.001269591250272884
$ws.Cells.Item(8, 13).Value = 22.618885
$ws.Cells.Item(8, 14).Value = 67.856655
$ws.Cells.Item(8, 15).Value = 0.9027998993061069
$ws.Cells.Item(8, 16).Value = 0.902799899306107
$ws.Cells.Item(8, 17).Value = 4.036777326693334
$ws.Cells.Item(8, 18).Value = 36.33099594024
$ws.Cells.Item(8, 19).Value = 0.001146186852906274
$ws.Cells.Item(8, 20).Value = 0.001146186852906274

# Row 9
$ws.Cells.Item(9, 7).Value = 0.1784693333333333
$ws.Cells.Item(9, 8).Value = 0.535408
$ws.Cells.Item(9, 9).Value = 0.001269591250272884
$ws.Cells.Item(9, 10).Value = 0.001269591250272884
$ws.Cells.Item(9, 13).Value = 2.154739666666666
$ws.Cells.Item(9, 14).Value = 6.464219
$ws.Cells.Item(9, 15).Value = 0.08600330007856447
$ws.Cells.Item(9, 16).Value = 0.08600330007856449
$ws.Cells.Item(9, 17).Value = 0.3845549518168889
$ws.Cells.Item(9, 18).Value = 3.460994566352
$ws.Cells.Item(9, 19).Value = 0.0001091890372743387
$ws.Cells.Item(9, 20).Value = 0.0001091890372743387

# Row 10
$ws.Cells.Item(10, 7).Value = 0.1784693333333333
$ws.Cells.Item(10, 8).Value = 0.535408
$ws.Cells.Item(10, 9).Value = 0.001269591250272884
$ws.Cells.Item(10, 10).Value = 0.001269591250272884
$ws.Cells.Item(10, 13).Value = 0.2805263333333333
$ws.Cells.Item(10, 14).Value = 0.841579
$ws.Cells.Item(10, 15).Value = 0.0111968006153285
$ws.Cells.Item(10, 16).Value = 0.01119680061532851
$ws.Cells.Item(10, 17).Value = 0.05006534769244445
$ws.Cells.Item(10, 18).Value = 0.450588129232
$ws.Cells.Item(10, 19).Value = 0.00001421536009227112
$ws.Cells.Item(10, 20).Value = 0.00001421536009227112

# New rows 11-13 for sending cluster "Resolving-Mac"
# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Sfrp1"
$ws.Cells.Item(11, 3).Value = "Fzd6"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.1385303333333333
$ws.Cells.Item(11, 8).Value = 0.415591
$ws.Cells.Item(11, 9).Value = 0.0009854740633164956
$ws.Cells.Item(11, 10).Value = 0.0009854740633164956
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 22.618885
$ws.Cells.Item(11, 14).Value = 67.856655
$ws.Cells.Item(11, 15).Value = 0.9027998993061069
$ws.Cells.Item(11, 16).Value = 0.902799899306107
$ws.Cells.Item(11, 17).Value = 3.133401678678334
$ws.Cells.Item(11, 18).Value = 28.200615108105
$ws.Cells.Item(11, 19).Value = 0.0008896858851309122
$ws.Cells.Item(11, 20).Value = 0.0008896858851309123

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Sfrp1"
$ws.Cells.Item(12, 3).Value = "Fzd6"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.1385303333333333
$ws.Cells.Item(12, 8).Value = 0.415591
$ws.Cells.Item(12, 9).Value = 0.0009854740633164956
$ws.Cells.Item(12, 10).Value = 0.0009854740633164956
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.154739666666666
$ws.Cells.Item(12, 14).Value = 6.464219
$ws.Cells.Item(12, 15).Value = 0.08600330007856447
$ws.Cells.Item(12, 16).Value = 0.08600330007856449
$ws.Cells.Item(12, 17).Value = 0.2984968042698889
$ws.Cells.Item(12, 18).Value = 2.686471238429
$ws.Cells.Item(12, 19).Value = 0.0000847540215870508
$ws.Cells.Item(12, 20).Value = 0.00008475402158705083

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Sfrp1"
$ws.Cells.Item(13, 3).Value = "Fzd6"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.1385303333333333
$ws.Cells.Item(13, 8).Value = 0.415591
$ws.Cells.Item(13, 9).Value = 0.0009854740633164956
$ws.Cells.Item(13, 10).Value = 0.0009854740633164956
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.2805263333333333
$ws.Cells.Item(13, 14).Value = 0.841579
$ws.Cells.Item(13, 15).Value = 0.0111968006153285
$ws.Cells.Item(13, 16).Value = 0.01119680061532851
$ws.Cells.Item(13, 17).Value = 0.03886140646544444
$ws.Cells.Item(13, 18).Value = 0.349752658189
$ws.Cells.Item(13, 19).Value = 0.00001103415659853242
$ws.Cells.Item(13, 20).Value = 0.00001103415659853242

